# Weekly update: a new "Poroto granado" price record (dated 2023-01-17,
# serial 44943) is inserted at the top of the existing data block (row 79),
# pushing the previously-existing rows 79-96 down to rows 80-97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 79 - shifts rows 79:96 down to 80:97 and
# grows the sheet dimension from A1:R96 to A1:R97.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Cells.Item(79, 1).Value = 10
$ws.Cells.Item(79, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(79, 3).Value = "La Araucanía"
$ws.Cells.Item(79, 4).Value = 44943
$ws.Cells.Item(79, 5).Value = 9
$ws.Cells.Item(79, 6).Value = 100112030
$ws.Cells.Item(79, 7).Value = "Poroto granado"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 65
$ws.Cells.Item(79, 11).Value = 45000
$ws.Cells.Item(79, 12).Value = 45000
$ws.Cells.Item(79, 13).Value = 45000
$ws.Cells.Item(79, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(79, 15).Value = "Región del Maule"
$ws.Cells.Item(79, 16).Value = 1800
$ws.Cells.Item(79, 17).Value = 25
$ws.Cells.Item(79, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
